$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Brasov/Hiking -> FakeCity/Dance
$ws.Range("B3").Value = "FakeCity"
$ws.Range("C3").Value = "Dance"

# Row 12 (Vaslui / Scuba Diving / 3) is removed entirely; remaining rows
# shift up, so former row 13 (Napoli / Pizza Making / 4.9) becomes row 12.
$ws.Rows("12").Delete()

# Rename the city on the (now) last row from "Napoli" to "Naples"
$ws.Range("B12").Value = "Naples"

# Update selection to match the new active cell
$ws.Range("F9").Select()
